$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "D61E3C1E"
$ws.Range("A3").Value = "BC964FD3"
$ws.Range("B3").Value = "CASSIANA UNICORNIO"
$ws.Range("C3").Value = 8

foreach ($r in @(5, 6, 7, 50, 51, 71)) {
    $ws.Range("A$r").Font.Bold = $false
}

$ws.Range("A3").Select()
